$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. Font change: Normal style font Arial -> Calibri (affects whole workbook)
# ---------------------------------------------------------------------------
$wb.Styles.Item("Normal").Font.Name = "Calibri"

# ---------------------------------------------------------------------------
# 2. Sheet1: add two new rows (8, 9 -> C10/C11) by copying the formatting of
#    the existing last data row (row 9) and then editing the values in the
#    same sequence the original author must have used so the shared-string
#    table is rebuilt in the same order as the target workbook.
# ---------------------------------------------------------------------------
$ws1.Range("C9:G9").Copy($ws1.Range("C10:G10"))
$ws1.Range("C9:G9").Copy($ws1.Range("C11:G11"))

$ws1.Range("C10").Value = 8
$ws1.Range("F10").Value = "getKhoaByMaBoPhan()"
$ws1.Range("D10").Value = "DAO"
$ws1.Range("D11").Value = "DAO"
$ws1.Range("E11").Value = "ThanhVienDAO"
$ws1.Range("F11").Value = "UpdateVaiTroTV()"
$ws1.Range("E10").Value = "KhoaDAO"
$ws1.Range("C11").Value = 9

# ---------------------------------------------------------------------------
# 3. Row heights: rows without an explicit customHeight re-flowed (font
#    metric change) from a 14.25 default row height to a 15 default row
#    height, i.e. scaled by 15/14.25.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(5).RowHeight = 135
$ws1.Rows.Item(6).RowHeight = 60
$ws1.Rows.Item(7).RowHeight = 120
$ws1.Rows.Item(8).RowHeight = 135

# ---------------------------------------------------------------------------
# 4. Column widths: small re-flow caused by the font-metric change.
# ---------------------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 4.714285714285714
$ws1.Columns.Item(3).ColumnWidth = 14.571428571428571
$ws1.Columns.Item(4).ColumnWidth = 27.714285714285715
$ws1.Columns.Item(5).ColumnWidth = 31.571428571428573
$ws1.Columns.Item(6).ColumnWidth = 61.142857142857146
$ws1.Columns.Item(7).ColumnWidth = 11.714285714285714

# ---------------------------------------------------------------------------
# 5. View state: window/scroll position.
# ---------------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.Width = 19440
